# Update clothes_2022-07-16.xlsx
#  - refresh the scrape timestamp (column O) for every data row
#  - a handful of products were re-scraped in a different order; their data
#    (columns A-N) needs to be moved to different rows (some rows simply
#    swap places, others rotate through 3 or 4 rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-07-16 20:56:48"

# ---------------------------------------------------------------------------
# Helpers to read/write the content of columns A:N for a given row while
# preserving data types:
#   - columns E and F hold numbers (ratingAmount / ratingValue)
#   - all the other columns (A-D, G-N) hold text; some of it looks numeric
#     (ids, prices) and must stay text, so the destination cell format is
#     forced to Text ("@") before writing the value and then restored to
#     whatever number format it had before.
# All reads for an operation are done up-front (into PowerShell variables)
# before any writes happen, so no temporary/helper cells are ever touched
# on the worksheet.
# ---------------------------------------------------------------------------
function Get-RowContent($row) {
    return @{
        Text1 = $ws.Range("A$row`:D$row").Value2
        Text2 = $ws.Range("G$row`:N$row").Value2
        Num   = $ws.Range("E$row`:F$row").Value2
    }
}

function Set-RowContent($row, $data) {
    $rngText1 = $ws.Range("A$row`:D$row")
    $rngText2 = $ws.Range("G$row`:N$row")
    $rngNum   = $ws.Range("E$row`:F$row")

    $fmt1 = $rngText1.NumberFormat
    $fmt2 = $rngText2.NumberFormat

    $rngText1.NumberFormat = "@"
    $rngText2.NumberFormat = "@"

    $rngText1.Value2 = $data.Text1
    $rngText2.Value2 = $data.Text2
    $rngNum.Value2   = $data.Num

    $rngText1.NumberFormat = $fmt1
    $rngText2.NumberFormat = $fmt2
}

# Rows 14 / 15 were re-scraped and ended up swapped
$d14 = Get-RowContent 14
$d15 = Get-RowContent 15
Set-RowContent 14 $d15
Set-RowContent 15 $d14

# Rows 22, 23, 24 rotated: new22 = old24, new23 = old22, new24 = old23
$d22 = Get-RowContent 22
$d23 = Get-RowContent 23
$d24 = Get-RowContent 24
Set-RowContent 22 $d24
Set-RowContent 23 $d22
Set-RowContent 24 $d23

# Rows 30, 31, 32, 33 rotated: new30 = old31, new31 = old32, new32 = old33, new33 = old30
$d30 = Get-RowContent 30
$d31 = Get-RowContent 31
$d32 = Get-RowContent 32
$d33 = Get-RowContent 33
Set-RowContent 30 $d31
Set-RowContent 31 $d32
Set-RowContent 32 $d33
Set-RowContent 33 $d30

# Rows 36 / 37 swapped
$d36 = Get-RowContent 36
$d37 = Get-RowContent 37
Set-RowContent 36 $d37
Set-RowContent 37 $d36

# Refresh the timestamp column for every data row (2-130)
$ws.Range("O2:O130").Value2 = $newTimestamp

Write-Host "Workbook updated"
